# Generate Report for Handoff
#
# Refreshes the "Latest Handoff Date(time)" values for the rows whose
# handoff finished just after the report was originally generated
# (the 1e39392e / 1afbb6e9 / 1ed63e7b / 3ea2fb8b / 455b92ba / 6ba422fa /
# b582b0a3 / e69672fa file group), moving their timestamp forward by a
# a minute or so, on the Overview sheet as well as the per-locale
# (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 7, 10, 11, 12, 13, 14, 15, 16) {
    $wsOverview.Range("D$r").Value = "2016-24-18 14:24:27"
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 7, 10, 11, 12, 13, 14, 15, 16) {
    $wsZhCn.Range("E$r").Value = "2016-03-18 14:24:23"
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in 7, 10, 11, 12, 13, 14, 15, 16) {
    $wsDeDe.Range("E$r").Value = "2016-03-18 14:24:27"
}
